# New weekly price record for "Femacal de La Calera - Ciboulette".
#
# A new observation (dated 2022-04-25, serial 44676) is inserted as row 63.
# Excel's Rows(...).Insert() shifts the existing rows 63:305 down to 64:306
# (pushing the sheet's used range from A1:R305 to A1:R306), carrying every
# cell's value/format along with it - which is exactly the "old row N -> new
# row N+1" pattern seen throughout the diff. The freshly inserted row 63 is
# then populated as a duplicate of the record that is now sitting in row 64
# (same market/region/category/quality/volume/prices/unit/origin/etc.), with
# only its date (column D) replaced by the new record's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 63
$newDateSerial = 44676

# Push rows 63:305 down to 64:306 (auto-extends the sheet to A1:R306).
$ws.Rows($newRow).Insert()

# The data that used to live in row 63 now lives in row $newRow + 1 - copy it
# across the full A:R record width to seed the newly inserted row.
$copyFromRow = $newRow + 1
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item($newRow, $col).Value2 = $ws.Cells.Item($copyFromRow, $col).Value2
}

# ... except the date, which is this new weekly record's own date (column D = 4).
$ws.Cells.Item($newRow, 4).Value2 = $newDateSerial
